$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number bump (si #6): "...Number  1" -> "...Number  2" ---
$ws.Range("A8").Characters(21, 1).Text = "2"

# --- Header: report-week date range (si #9) ---
# Replace right-hand date first so the left offset (27) stays valid.
$ws.Range("C9").Characters(48, 8).Text = "1/12/2025"
$ws.Range("C9").Characters(27, 10).Text = "1/6/2025"

# --- CompStat table rows 15-28: Week-to-Date / 28-Day / YTD / 2-Year figures ---
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("N15").Value = 0
$ws.Range("H15").Copy()
$ws.Range("N15").PasteSpecial(-4122)

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 6
$ws.Range("F15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -83.333333333333
$ws.Range("H15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 6
$ws.Range("F15").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$ws.Range("K16").Value = -33.333333333333
$ws.Range("H15").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = 300
$ws.Range("H15").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("N16").Value = -86.206896551724

$ws.Range("C17").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -30
$ws.Range("M17").Value = -100
$ws.Range("H15").Copy()
$ws.Range("M17").PasteSpecial(-4122)

$ws.Range("C18").Value = 4
$ws.Range("F15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 4
$ws.Range("F15").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = -20
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = -66.666666666666
$ws.Range("N18").Value = -88.235294117647

$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = -10.416666666666
$ws.Range("I19").Value = 32
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = -13.513513513513
$ws.Range("L19").Value = -3.030303030303
$ws.Range("M19").Value = 6.666666666666
$ws.Range("N19").Value = -70.909090909090

$ws.Range("D20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -100

$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -7.407407407407
$ws.Range("F21").Value = 116
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = -16.546762589928
$ws.Range("I21").Value = 41
$ws.Range("J21").Value = 51
$ws.Range("K21").Value = -19.607843137254
$ws.Range("L21").Value = -16.326530612244
$ws.Range("M21").Value = -8.888888888888
$ws.Range("N21").Value = -80.382775119617

$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -20
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = -66.666666666666
$ws.Range("L22").Value = -80
$ws.Range("M22").Value = -50
$ws.Range("H15").Copy()
$ws.Range("M22").PasteSpecial(-4122)

$ws.Range("C24").Value = 69
$ws.Range("D24").Value = 92
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 274
$ws.Range("G24").Value = 329
$ws.Range("H24").Value = -16.717325227963
$ws.Range("I24").Value = 105
$ws.Range("J24").Value = 145
$ws.Range("K24").Value = -27.586206896551
$ws.Range("L24").Value = -18.604651162790
$ws.Range("M24").Value = 114.285714285714

$ws.Range("C25").Value = 70
$ws.Range("D25").Value = 90
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 270
$ws.Range("G25").Value = 341
$ws.Range("H25").Value = -20.821114369501
$ws.Range("I25").Value = 104
$ws.Range("J25").Value = 146
$ws.Range("K25").Value = -28.767123287671
$ws.Range("L25").Value = -18.110236220472

$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 250
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 11.538461538461
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -20
$ws.Range("M26").Value = 100

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 1
$ws.Range("F15").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("J28").Value = 2
$ws.Range("F15").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("K28").Value = -50
$ws.Range("H15").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("L28").Value = -50
